$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: timestamp text update ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 17:35"

# --- Estados Unidos (row 4): updated case counts ---
$ws.Range("B4").Value = 1670402
$ws.Range("C4").Value = 3574
$ws.Range("E4").Value = 1124668
$ws.Range("G4").Value = 67
$ws.Range("H4").Value = 98750

# --- Canada (row 16): updated case counts ---
$ws.Range("B16").Value = 84081
$ws.Range("C16").Value = 460
$ws.Range("E16").Value = 34396
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 6380

# --- Rows 124-125: Sierra Leona / Jordania swap places ---
# Row 124 becomes Jordania with freshly updated data; row 125 becomes
# Sierra Leona carrying the data previously on row 124 unchanged.
$ws.Range("A124").Value = "Jordania"
$ws.Range("B124").Value = 708
$ws.Range("C124").Value = 4
$ws.Range("D124").Value = 471
$ws.Range("E124").Value = 228
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 9

$ws.Range("A125").Value = "Sierra Leona"
$ws.Range("B125").Value = 707
$ws.Range("C125").Value = 86
$ws.Range("D125").Value = 241
$ws.Range("E125").Value = 426
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 40

# --- Rows 198-200: Nueva Caledonia / Belice / Santa Lucia rotate ---
$ws.Range("A198").Value = "Belice"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 16
$ws.Range("E198").Value = 0
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 2

$ws.Range("A199").Value = "Santa Lucia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# --- Rows 209-211: Seychelles / Groenlandia / Montserrat rotate ---
$ws.Range("A209").Value = "Montserrat"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 10
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# --- Rows 214-216: Bonaire.. / San Bartolome / Sahara Occidental rotate ---
# (numeric data identical across the three, only names move)
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
